# Fruta / hortaliza, semanal
# Inserts two new weekly price records (Chirimoya, Macroferia Regional de
# Talca) into the middle of the existing table, shifting the subsequent
# rows down - matching the target diff exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new record at row 11 --------------------------------------
# (pushes the former row 11 and everything below it down by one row)
$ws.Rows.Item(11).Insert()

$ws.Cells.Item(11,1).Value2  = 5
$ws.Cells.Item(11,2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(11,3).Value2  = "Maule"
$ws.Cells.Item(11,4).Value2  = 44482
$ws.Cells.Item(11,5).Value2  = 7
$ws.Cells.Item(11,6).Value2  = "Fruta"
$ws.Cells.Item(11,7).Value2  = 100107
$ws.Cells.Item(11,8).Value2  = "Otros"
$ws.Cells.Item(11,9).Value2  = 100107002
$ws.Cells.Item(11,10).Value2 = "Chirimoya"
$ws.Cells.Item(11,11).Value2 = "Cultivar IV Región"
$ws.Cells.Item(11,12).Value2 = "Primera"
$ws.Cells.Item(11,13).Value2 = 120
$ws.Cells.Item(11,14).Value2 = 25000
$ws.Cells.Item(11,15).Value2 = 25000
$ws.Cells.Item(11,16).Value2 = 25000
$ws.Cells.Item(11,17).Value2 = '$/bandeja 10 kilos'
$ws.Cells.Item(11,18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(11,19).Value2 = 2500
$ws.Cells.Item(11,20).Value2 = 10

# --- Insert new record at row 15 --------------------------------------
# (pushes the former row 14 [already shifted once, now at 14] and
# everything below it down by one more row)
$ws.Rows.Item(15).Insert()

$ws.Cells.Item(15,1).Value2  = 5
$ws.Cells.Item(15,2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(15,3).Value2  = "Maule"
$ws.Cells.Item(15,4).Value2  = 44421
$ws.Cells.Item(15,5).Value2  = 7
$ws.Cells.Item(15,6).Value2  = "Fruta"
$ws.Cells.Item(15,7).Value2  = 100107
$ws.Cells.Item(15,8).Value2  = "Otros"
$ws.Cells.Item(15,9).Value2  = 100107002
$ws.Cells.Item(15,10).Value2 = "Chirimoya"
$ws.Cells.Item(15,11).Value2 = "Cultivar IV Región"
$ws.Cells.Item(15,12).Value2 = "Especial"
$ws.Cells.Item(15,13).Value2 = 30
$ws.Cells.Item(15,14).Value2 = 35000
$ws.Cells.Item(15,15).Value2 = 35000
$ws.Cells.Item(15,16).Value2 = 35000
$ws.Cells.Item(15,17).Value2 = '$/bandeja 10 kilos'
$ws.Cells.Item(15,18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(15,19).Value2 = 3500
$ws.Cells.Item(15,20).Value2 = 10
